$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = '宇宙兄弟猫の口から煙をやりだしとやります しかしんこやりますやりますので 各自 ヘッド 3万円で売りましょう'
$ws.Range('C3').Value = '自分から飯します まず このスパイダーマンの鉛筆を3万円で売る方法として まず リアル スパイダーマンの鉛筆を持っこの行ってたりして めちゃくちゃ 付加価値を高めて3万円で売るというほうがいいと思います できれば 鉛筆の価値を知らない人に売りつけるのが一番手っ取り早いと思って あ このスパイダーマン この鉛筆の価値を決めるのって手間 スパイダーマンというブランドがあるんでめちゃくちゃ 昔のあのコレクターが集めるような12 希少性の高い鉛筆という説明をして価値を知らない人に'
$ws.Range('C4').Value = '価値を知らない人に3万円で売りつけるという方法を考えましたか'
$ws.Range('C5').Value = '平均額になるも俺はそれで 智子 の価値をつけてるっていう空 トムホランドが使った1ドルって この後なんかなんか'
$ws.Range('C6').Value = '今時間どれくらいですか 10秒まで行くよ 5分 あとからまたこのパイダーマンの鉛筆の消しゴムと鉛筆の周りを覆っている この金属部分が 実は銀でできてるとか'
$ws.Range('C7').Value = 'そんな誰かが入ってるから中 日本語が限られてたりとかしてる可能性とかもあるんじゃないかなと思ってやっぱ3万円で売るのが一番いいのではないかと基本的には茶色い 鉛筆がデフォルトの鉛筆やからそれはまき'
$ws.Range('C8').Value = '鉛筆やからそれは基本的に そんな高くはないと思うんですけど このスパイダーマンのこの鉛筆は売ってる場所とかもうそんな前全国でそんないろいろいろんな場所 取ってるわけではないからもう希少価値が高くなってるって言う事で3万で売るべきではないでしょうか 絵柄だけでさんま 買ってきて今売ってる場所も限られてるのもあるからです 体で全国で'
$ws.Range('C9').Value = '全国でも合う文房具屋さんとかなんかたくさんところで 市販で売ってるけれども スパイダー 真夏はその数店舗しか売ってないことによって価値がちょっとあるんではないかと思いますが そこまで手が回らん状態で高い 付加価値つけて'
$ws.Range('C10').Value = 'しましょう も 怒んないのが そのメルカリでは何かしら 松前 太田君が入った その方が使用してたとかで高値でまたますっていうのは あれやろけどそれでも メルカリでそういう 紹介文に書いて売るなら 一番なんか売れるのかな っていうのしか 俺は愛しかない前に一通りの意見でたんですけど 7分も出してみる一応 回答は 企業側の提案やから真っ当な理由で見つけないとやばいは多分10倍とか言ったら怒られそう'
$ws.Range('C11').Value = 'お前の意見をまとめる そのまあどういう風にしてこの鉛筆に今 価値をつける かすればいいし お菓子だったりまーす 形だったりをつけることによっては3万円で売るっていうふうにしてると思うんだけども 一番 現実的なのはさっきあったようにその俳優さんとかが使ったもんとかなんか サイン書いてます とかそういうので3万円で売る での 噛まれた後の処置は懸念点としてはないあの そのまま買う相手の人が誰かしら 俳優 だとしても パンパンでなかったり興味がなかったりしたら あんまんだよ 買わんと思うよ その人になって 今勝手が違うって言う'
$ws.Range('C12').Value = '違うって言うとそうそうそう それからやっぱり最初に鉛筆そのもの自体を これ これだけにするっていうのはありだと思うけどそれだったら一応その案とか関係ないっていう話いるんやったらさっき この鉛筆 自体に 銀と香は鉱物が使われてると帰ったらもうなんかそういうファンとか スパイダーマンとか関係なしにもう一つしてついてるものと想定されるから その売掛 手帳 3万 行けるんじゃないか そのお金自体で3万円するみたいな'
$ws.Range('C13').Value = '今日は カテキンでなさそうで 止めようかなと 今でもそれって それぞれの良いところをまとめるとかでもないもんな これ そうじゃないと現実としては すべての人にダメって言うか 家で売るなら また金とか半分にして'
$ws.Range('C14').Value = 'マイミクにしてるって言うのがまた後じゃないかなと思います どうしたんですので もしやな ただ 実際にでも銀でできてるわけじゃない なんでの問題の受け取り方としたらさ お前はぎんじゃない鉛筆を3万円で売ってください っていう意味やったとしたらこういう事情もあった2人 正直この鉛筆を通したんですけど この鉛筆を3万円でお金を持ってる人に売りつけてで売ったお金は'
$ws.Range('C15').Value = 'で売ったお金は慈善活動に使われますという名目のもと 3万円で売る 競ったらちゃん前に買ってくれんじゃなかったらお金を事前活を入れるだけじゃなくてお金持ってる人がもう スパイダーマンえんぴつが満を持してて入るわけやしてない子全体的にもお金が欲しい 支援が言う人とかにも使われるお金がでビンビンな関係 やんか こら 鉛筆を3万円で買っていただけると 慈善活動につながりますって言う 気をつけて 売りつけるっていうのが今思い浮かんだんじゃないか 時間ないけどあんまりあと1分からもう一個なんか考えて 鉛筆に何かしら 他の なんか サービスを受けて'
$ws.Range('C16').Value = '鉛筆 自体の値段 3万 やけど そのサービスに入ってもらったら3万キャッシュバックっていう形で入ってくる うるうるって言うのだったら 鉛筆 実質タダ みたいなもんやけど満々で売ることはできるけどちょっと まとめてきても理想 何がいいかなっていうのはこの鉛筆の3万円に対して何かしらのコンテンツやサービスを付加することによってその鉛筆に付属するものを3枚の価値として付与することで3万円 出る 鉛筆を売りつけるという形でいいんじゃないかな って気になりました います'

$ws.Rows("17:20").Delete()

Write-Output $ws.UsedRange.Address()
